$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 30
$ws.Range("E6").Value = 33
$ws.Range("F6").Value = 20
$ws.Range("H6").Value = 20
$ws.Range("E10").Value = 302
$ws.Range("F10").Value = 151
$ws.Range("H10").Value = 151
$ws.Range("E11").Value = 210
$ws.Range("F11").Value = 126
$ws.Range("H11").Value = 126
$ws.Range("E12").Value = 311
$ws.Range("F12").Value = 185
$ws.Range("H12").Value = 185
$ws.Range("E13").Value = 93
$ws.Range("F13").Value = 49
$ws.Range("H13").Value = 49
$ws.Range("E14").Value = 87
$ws.Range("E15").Value = 111
$ws.Range("F15").Value = 45
$ws.Range("H15").Value = 45
$ws.Range("E16").Value = 129
$ws.Range("F16").Value = 74
$ws.Range("H16").Value = 74
$ws.Range("E17").Value = 58
$ws.Range("F17").Value = 31
$ws.Range("H17").Value = 31
$ws.Range("E20").Value = 66
$ws.Range("E21").Value = 96
$ws.Range("F21").Value = 55
$ws.Range("H21").Value = 55
$ws.Range("E22").Value = 116
$ws.Range("F22").Value = 69
$ws.Range("H22").Value = 69
$ws.Range("E23").Value = 121
$ws.Range("F23").Value = 65
$ws.Range("H23").Value = 65
$ws.Range("F24").Value = 73
$ws.Range("H24").Value = 73
$ws.Range("E25").Value = 153
$ws.Range("F25").Value = 77
$ws.Range("H25").Value = 77
$ws.Range("F26").Value = 55
$ws.Range("H26").Value = 55
$ws.Range("E27").Value = 200
$ws.Range("E28").Value = 121
$ws.Range("E29").Value = 121
$ws.Range("F30").Value = 82
$ws.Range("H30").Value = 82
$ws.Range("E31").Value = 52
$ws.Range("F32").Value = 71
$ws.Range("H32").Value = 71
$ws.Range("E33").Value = 192
$ws.Range("F33").Value = 102
$ws.Range("H33").Value = 102
$ws.Range("E34").Value = 144
$ws.Range("F34").Value = 92
$ws.Range("H34").Value = 92
$ws.Range("F36").Value = 29
$ws.Range("H36").Value = 29
$ws.Range("E37").Value = 104
$ws.Range("F37").Value = 57
$ws.Range("H37").Value = 57
$ws.Range("E38").Value = 63
$ws.Range("E39").Value = 126
$ws.Range("E40").Value = 179
$ws.Range("E41").Value = 247
$ws.Range("F41").Value = 116
$ws.Range("H41").Value = 116
$ws.Range("E42").Value = 230
$ws.Range("F42").Value = 125
$ws.Range("H42").Value = 125
$ws.Range("F43").Value = 39
$ws.Range("H43").Value = 39
$ws.Range("E44").Value = 196
$ws.Range("F44").Value = 110
$ws.Range("H44").Value = 110
$ws.Range("E45").Value = 76
$ws.Range("F45").Value = 44
$ws.Range("H45").Value = 44
$ws.Range("E47").Value = 291
$ws.Range("F47").Value = 149
$ws.Range("H47").Value = 149
$ws.Range("E49").Value = 167
$ws.Range("F49").Value = 86
$ws.Range("H49").Value = 86
$ws.Range("E50").Value = 144
$ws.Range("F50").Value = 59
$ws.Range("H50").Value = 59
$ws.Range("E51").Value = 139
$ws.Range("F51").Value = 64
$ws.Range("H51").Value = 64
